$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 8000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H113").Value = 3422.5
$ws.Range("I113").Value = 3405
$ws.Range("J113").Value = 3475
$ws.Range("K113").Value = 3405
$ws.Range("L113").Value = 3475
$ws.Range("M113").Value = -151
$ws.Range("N113").Value = -9983

$ws.Range("H129").Value = 993.3333
$ws.Range("I129").Value = 636.6667
$ws.Range("J129").Value = 1025.7576
$ws.Range("K129").Value = 1910.0001
$ws.Range("L129").Value = 3077.2728
$ws.Range("M129").Value = 3089.9999
$ws.Range("N129").Value = -13077.2728

$ws.Range("H137").Value = 3451132
$ws.Range("I137").Value = 3848185.8
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 11544557.4
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -11542007.4
$ws.Range("N137").Value = -35100

$ws.Range("H138").Value = 1963248.4
$ws.Range("I138").Value = 691.71875
$ws.Range("J138").Value = 3148188.2
$ws.Range("K138").Value = 2075.15625
$ws.Range("L138").Value = 9444564.600000001
$ws.Range("M138").Value = 3064.84375
$ws.Range("N138").Value = -9454844.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2840.7144
$ws.Range("I63").Value = 2840.7144
$ws.Range("K63").Value = 2840.7144
$ws.Range("M63").Value = -2154.7144

$ws.Range("H66").Value = 2840.7144
$ws.Range("I66").Value = 2840.7144
$ws.Range("K66").Value = 14203.572
$ws.Range("M66").Value = -10771.572

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1065.2222
$ws.Range("I99").Value = 1063.762
$ws.Range("K99").Value = 1063.762
$ws.Range("M99").Value = 434.2380000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 650
$ws.Range("I2").Value = 650
$ws.Range("K2").Value = 650
$ws.Range("M2").Value = -537

$ws.Range("H31").Value = 2427.4807
$ws.Range("I31").Value = 1460.9474
$ws.Range("J31").Value = 5050.9287
$ws.Range("K31").Value = 1460.9474
$ws.Range("L31").Value = 5050.9287
$ws.Range("M31").Value = -1165.9474
$ws.Range("N31").Value = -5640.9287

$ws.Range("H34").Value = 2427.4807
$ws.Range("I34").Value = 1460.9474
$ws.Range("J34").Value = 5050.9287
$ws.Range("K34").Value = 1460.9474
$ws.Range("L34").Value = 5050.9287
$ws.Range("M34").Value = -1258.9474
$ws.Range("N34").Value = -5454.9287

$ws.Range("H99").Value = 4221.8423
$ws.Range("I99").Value = 3540.4167
$ws.Range("J99").Value = 5390
$ws.Range("K99").Value = 3540.4167
$ws.Range("L99").Value = 5390
$ws.Range("M99").Value = -2042.4167
$ws.Range("N99").Value = -8386

$ws.Range("H126").Value = 4221.8423
$ws.Range("I126").Value = 3540.4167
$ws.Range("J126").Value = 5390
$ws.Range("K126").Value = 10621.2501
$ws.Range("L126").Value = 16170
$ws.Range("M126").Value = -8151.250100000001
$ws.Range("N126").Value = -21110

$ws.Range("H132").Value = 14966.361
$ws.Range("I132").Value = 936.5246
$ws.Range("J132").Value = 92768.17999999999
$ws.Range("K132").Value = 2809.5738
$ws.Range("L132").Value = 278304.54
$ws.Range("M132").Value = -279.5738000000001
$ws.Range("N132").Value = -283364.54

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 543.1818
$ws.Range("I6").Value = 53.57143
$ws.Range("K6").Value = 160.71429
$ws.Range("M6").Value = -47.71429000000001

$ws.Range("H96").Value = 4600
$ws.Range("J96").Value = 4600
$ws.Range("L96").Value = 13800
$ws.Range("N96").Value = -17918

$ws.Range("H100").Value = 2666.3333
$ws.Range("I100").Value = 1800
$ws.Range("J100").Value = 2745.0908
$ws.Range("K100").Value = 5400
$ws.Range("L100").Value = 8235.2724
$ws.Range("M100").Value = -4589
$ws.Range("N100").Value = -9857.2724

$ws.Range("H106").Value = 3960.2
$ws.Range("J106").Value = 3960.2
$ws.Range("L106").Value = 11880.6
$ws.Range("N106").Value = -13772.6

$ws.Range("H109").Value = 3158.25
$ws.Range("I109").Value = 2616.5
$ws.Range("J109").Value = 3700
$ws.Range("K109").Value = 7849.5
$ws.Range("L109").Value = 11100
$ws.Range("M109").Value = -6809.5
$ws.Range("N109").Value = -13180

$ws.Range("H112").Value = 23811990
$ws.Range("I112").Value = 2129.25
$ws.Range("J112").Value = 55558470
$ws.Range("K112").Value = 6387.75
$ws.Range("L112").Value = 166675410
$ws.Range("M112").Value = -5279.75
$ws.Range("N112").Value = -166677626

$ws.Range("H121").Value = 58850300
$ws.Range("I121").Value = 1260
$ws.Range("J121").Value = 70240430
$ws.Range("K121").Value = 3780
$ws.Range("L121").Value = 210721290
$ws.Range("M121").Value = -2470
$ws.Range("N121").Value = -210723910

$ws.Range("H131").Value = 1381
$ws.Range("I131").Value = 775
$ws.Range("J131").Value = 1424.2858
$ws.Range("K131").Value = 2325
$ws.Range("L131").Value = 4272.857400000001
$ws.Range("M131").Value = 2715
$ws.Range("N131").Value = -14352.8574

$ws.Range("H132").Value = 1324.1428
$ws.Range("I132").Value = 696.6667
$ws.Range("J132").Value = 1794.75
$ws.Range("K132").Value = 6270.0003
$ws.Range("L132").Value = 16152.75
$ws.Range("M132").Value = -3740.0003
$ws.Range("N132").Value = -21212.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2656.3784
$ws.Range("I122").Value = 2302.5
$ws.Range("J122").Value = 3757.3333
$ws.Range("K122").Value = 6907.5
$ws.Range("L122").Value = 11271.9999
$ws.Range("M122").Value = -4457.5
$ws.Range("N122").Value = -16171.9999

$ws.Range("H126").Value = 3578.25
$ws.Range("I126").Value = 5204
$ws.Range("J126").Value = 2602.8
$ws.Range("K126").Value = 15612
$ws.Range("L126").Value = 7808.400000000001
$ws.Range("M126").Value = -13142
$ws.Range("N126").Value = -12748.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1471.5
$ws.Range("I100").Value = 1193.4667
$ws.Range("J100").Value = 1792.3077
$ws.Range("K100").Value = 1193.4667
$ws.Range("L100").Value = 1792.3077
$ws.Range("M100").Value = -652.4666999999999
$ws.Range("N100").Value = -2874.3077

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1853.3334
$ws.Range("I122").Value = 1405.0358
$ws.Range("J122").Value = 2749.9285
$ws.Range("K122").Value = 4215.107400000001
$ws.Range("L122").Value = 8249.7855
$ws.Range("M122").Value = -1765.107400000001
$ws.Range("N122").Value = -13149.7855

$ws.Range("H126").Value = 1508.3636
$ws.Range("I126").Value = 1356
$ws.Range("J126").Value = 1775
$ws.Range("K126").Value = 4068
$ws.Range("L126").Value = 5325
$ws.Range("M126").Value = -1598
$ws.Range("N126").Value = -10265

$ws.Range("H136").Value = 33278.777
$ws.Range("I136").Value = 21155.26
$ws.Range("K136").Value = 63465.78
$ws.Range("M136").Value = -60915.78

